# Commit: "Day 4 not Day 2"
#
# The workbook had sheets "D2-Solution"/"D2-Input" that were actually the
# Day 4 puzzle (mislabeled as Day 2). Rename them to "D4-Solution"/
# "D4-Input" (Excel auto-updates the cross-sheet formula references that
# point at 'D2-Input' so they become 'D4-Input'), fix the "Day 2" header
# text on both sheets to "Day 4", and move the active sheet/selection from
# D2-Solution (now D4-Solution) to D2-Input (now D4-Input), matching the
# new cursor position saved in the file.

$wb = $excel.ActiveWorkbook

$wsSolution = $wb.Worksheets.Item("D2-Solution")
$wsInput = $wb.Worksheets.Item("D2-Input")

# Renaming updates every formula that referenced the old sheet names.
$wsSolution.Name = "D4-Solution"
$wsInput.Name = "D4-Input"

# Fix the visible "Day 2" label on both sheets.
$wsSolution.Range("A1").Value = "Day 4"
$wsInput.Range("A1").Value = "Day 4"

# D4-Solution keeps a plain selection at A2 (no longer the active tab).
$wsSolution.Activate()
$wsSolution.Range("A2").Select()

# D4-Input becomes the active tab, selection reset to A2.
$wsInput.Activate()
$wsInput.Range("A2").Select()
